$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the two new rows (2021年, 2022年), columns B..AQ (42 values each)
$row11 = @(264.54,51.77,13.18,0.08,356.4,825.7,121.23,193.19,51.7,56.27,59.95,4.91,92.37,566.95,64.48,65.7,219.83,201.38,100.8,-85.31,898.34,105.71,2968.29,452.49,64.87,380.33,163.48,266.67,100.42,11428.25,431.75,262.84,161.28,48.32,14.16,339.29,62.34,582.79,70.34,97.84999999999999,555.41,176.17)
$row12 = @(133.3,18.9,7.9,0,367,701.5,114.8,46,38.1,36.8,61.5,-2.8,57.6,473.8,65.09999999999999,59.1,149.4,234.8,83.90000000000001,-92.59999999999999,775.5,126.1,3148.9,131.1,47.2,421.1,151.1,225.7,48.8,9850.299999999999,198.9,146.6,151.7,34.1,8.300000000000001,271.3,32.4,550.7,75.09999999999999,63.5,536.8,151.1)

$years = @("2021年", "2022年")
$dataRows = @($row11, $row12)

for ($r = 0; $r -lt 2; $r++) {
    $targetRow = 11 + $r

    # Copy the formatting (font/border/alignment) of the year cell from the row above
    $ws.Range("A" + ($targetRow - 1)).Copy()
    $ws.Cells.Item($targetRow, 1).PasteSpecial(-4122)

    # Year label in column A
    $ws.Cells.Item($targetRow, 1).Value = $years[$r]

    # Data values in columns B (2) .. AQ (43)
    $vals = $dataRows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($targetRow, 2 + $i).Value = $vals[$i]
    }
}

$excel.CutCopyMode = 0
